$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Newly-added Installed Capacity ")

# --- Row 6 ("Time Period"): shift the YTD start:end window forward one month (…:2023-05 -> …:2023-06) ---
$ws.Range("B6").Value2 = "2008-12:2023-06"
$ws.Range("C6").Value2 = "2008-12:2023-06"
$ws.Range("D6").Value2 = "2010-09:2023-06"
$ws.Range("E6").Value2 = "2008-12:2023-06"
$ws.Range("F6").Value2 = "2009-12:2023-06"
$ws.Range("I6").Value2 = "2012-12:2023-06"
$ws.Range("J6").Value2 = "2012-12:2023-06"
$ws.Range("K6").Value2 = "2012-12:2023-06"
$ws.Range("L6").Value2 = "2012-12:2023-06"
$ws.Range("M6").Value2 = "2012-12:2023-06"
# G6,H6,N6,O6,P6,Q6 (Wind/Solar series) keep their prior "Time Period" value unchanged

# --- Row 8 ("Update"): most series refreshed 2023-06-21 -> 2023-07-21 ---
# (leading apostrophe keeps these stored as text, matching the source file,
#  instead of being auto-converted to a date serial)
$ws.Range("B8").Value2 = "'2023-07-21"
$ws.Range("C8").Value2 = "'2023-07-21"
$ws.Range("D8").Value2 = "'2023-07-21"
$ws.Range("E8").Value2 = "'2023-07-21"
$ws.Range("F8").Value2 = "'2023-07-21"
$ws.Range("I8").Value2 = "'2023-07-21"
$ws.Range("J8").Value2 = "'2023-07-21"
$ws.Range("K8").Value2 = "'2023-07-21"
$ws.Range("L8").Value2 = "'2023-07-21"
$ws.Range("M8").Value2 = "'2023-07-21"
# G8,H8,N8,O8,P8,Q8 (Wind/Solar series) keep their prior "Update" value unchanged

# --- Append the new monthly data row (2023-06 / serial 45107) ---
$ws.Range("A165:Q165").Copy()
$ws.Range("A166:Q166").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A166").Value2 = 45107
$newRow = @(536, 2602, 119, 2299, 7842, 0, 0, 270772, 41793, 135698, 5676, 38921, 0, 0, 0, 0)
$col = 2
foreach ($v in $newRow) {
  $ws.Cells.Item(166, $col).Value2 = $v
  $col = $col + 1
}

# --- Update the saved selection to match the source workbook ---
$ws.Range("F2").Select() | Out-Null
